$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New French-locale rows (10013-10018), mirroring the existing
# eng/ara "Pre-Registration / Registration Client / Registration
# Processor / ID Authentication / ID Control / Resident Portal" rows.
$newRows = @(
    @{ Row = 14; Id = 10013; Name = "Pré-inscription";          Descr = "Portail Web pour les pré-inscriptions" },
    @{ Row = 15; Id = 10014; Name = "Client dinscription";      Descr = "Application de bureau pour les inscriptions" },
    @{ Row = 16; Id = 10015; Name = "Processeur dinscription";  Descr = "Demande de post-inscription" },
    @{ Row = 17; Id = 10016; Name = "Authentification ID";      Descr = "Application pour lauthentification du fournisseur de services tiers" },
    @{ Row = 18; Id = 10017; Name = "Contrôle didentité";       Descr = "Portail Web pour la configuration dapplications" },
    @{ Row = 19; Id = 10018; Name = "Portail Résident";         Descr = "Portail Web pour les services de génération de post-ID" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Descr
    $ws.Cells.Item($row, 4).Value = "fra"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Column sizing (id column auto/best-fit-ish, name column widened)
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 18.5

# Scroll the view down and select the rows below the data, matching
# the author's on-save selection state.
[void]$ws.Rows("20:1048576").Select()

# Page setup: A4, portrait.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
